# Weekly update of the "Macroferia Regional de Talca - Zanahoria" sheet.
# A new observation is inserted at the top of the recent-history block
# (row 460) and every subsequent row (461..498) is shifted down by one,
# pushing the oldest row (old 498) out into a brand-new row 499.
#
# Columns that carry the actual weekly record (D=Fecha, J=Volumen,
# K=Precio minimo, L=Precio maximo, M=Precio promedio ponderado,
# O=Origen, P=Precio $/Kg) are the ones that move; the rest of each row
# (A,B,C,E,F,G,H,I,N,Q,R) is identical across the whole block, so it is
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: append the new last row (499) as an exact copy of the current
# last row (498), before anything else is modified.
$ws.Range("A498:R498").Copy($ws.Range("A499:R499"))

# Step 2: shift rows 461..498 down by one. Walk from the bottom (498) up
# to 461 so that, for each row r, row r-1 still holds its original
# (pre-shift) values when it is copied into row r.
for ($r = 498; $r -ge 461; $r--) {
    $prev = $r - 1

    $srcD  = "D" + $prev
    $dstD  = "D" + $r
    $srcJM = "J" + $prev + ":M" + $prev
    $dstJM = "J" + $r + ":M" + $r
    $srcOP = "O" + $prev + ":P" + $prev
    $dstOP = "O" + $r + ":P" + $r

    $ws.Range($srcD).Copy($ws.Range($dstD))
    $ws.Range($srcJM).Copy($ws.Range($dstJM))
    $ws.Range($srcOP).Copy($ws.Range($dstOP))
}

# Step 3: write the brand-new observation into row 460. Its origin
# (O460, "Region de Nuble") stays the same as before.
$ws.Range("D460").Value = 45013
$ws.Range("J460").Value = 500
$ws.Range("K460").Value = 7000
$ws.Range("L460").Value = 7000
$ws.Range("M460").Value = 7000
$ws.Range("P460").Value = 350
